$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename / replace the header row for the expense-import/export sheet.
$ws.Range("A1").Value = "expense_type"
$ws.Range("B1").Value = "store"
$ws.Range("C1").Value = "employee_code"
$ws.Range("D1").Value = "amount"
$ws.Range("E1").Value = "expense_date"
$ws.Range("F1").Value = "note"

# Mark the original three header cells with an explicit (re-applied) style
# so they get their own stamped cellXf, distinguishing them from the three
# newly-appended columns which keep the sheet's default formatting.
$ws.Range("A1").Style = "Normal"
$ws.Range("B1:C1").Style = "Normal"

# Widen the columns to fit the new header labels.
$ws.Columns.Item(1).ColumnWidth = 11.83
$ws.Columns.Item(2).ColumnWidth = 4.64
$ws.Columns.Item(3).ColumnWidth = 13.29
$ws.Columns.Item(4).ColumnWidth = 6.65
$ws.Columns.Item(5).ColumnWidth = 11.93
$ws.Columns.Item(6).ColumnWidth = 4.1

# Leave the selection where the author's last edit left it.
$ws.Range("A6").Select() | Out-Null
